# Revert "fix broken mapping_orig.txt links"
# (reverts commit 5378006d9595cfc87e9b5edafb1224d4c824809c)
#
# The workbook currently has the R-column "mapping-orig.txt" hyperlink
# formulas pointing at "../../blob/master/datasets/..." GitHub blob links
# (one standalone CONCATENATE() per row). The commit being reverted had
# "fixed" those into root-relative "./datasets/..." paths and, because the
# per-row formulas became byte-identical once the row-relative B-column ref
# is factored out, Excel collapsed R3:R16 into a shared-formula group
# (R2 stayed standalone). Reproduce that prior state: rewrite the formulas
# back to "./datasets/..." (still landing as shared R3:R16 + standalone R2)
# plus the view-state changes (active sheet/selection) that were saved
# alongside it.

$wb = $excel.ActiveWorkbook

$wsData  = $wb.Worksheets.Item("dataset_metadata")
$wsTasks = $wb.Worksheets.Item("tasks")

# --- Formulas: R2 (standalone) + R3:R16 (shared formula group) -------------
# Writing a relative formula across a contiguous range makes Excel
# auto-adjust the B-column reference per row and collapse the identical
# formulas into one shared-formula group (si=0) anchored at R3 -- matching
# the target XML exactly (R2 keeps its own <f>, R3:R16 share si="0").
$wsData.Range("R2").Formula = '=CONCATENATE("./datasets/", B2, "/mapping-orig.txt")'
$wsData.Range("R3:R16").Formula = '=CONCATENATE("./datasets/", B3, "/mapping-orig.txt")'

# --- View state -------------------------------------------------------------
# Selection moves to F22 on dataset_metadata, and "tasks" becomes the
# selected/active tab (dataset_metadata loses tabSelected, workbook's
# activeTab becomes 1 = "tasks").
$wsData.Range("F22").Select()
$wsTasks.Activate()

$wb.Save()
